$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Range("B1").Value = "ID FLEX"
